$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 4.9914799125796581
$ws.Range("C2").Value = 10.466253638017726
$ws.Range("D2").Value = 12.167178286248378
$ws.Range("E2").Value = 10.663079852511585

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 3.8442887910512433
$ws.Range("C3").Value = 6.0993965164398682
$ws.Range("D3").Value = 15.890322102211959
$ws.Range("E3").Value = 7.2841762501876959

# Update the selected range to match the new selection noted in the diff (B1:E3)
$ws.Range("B1:E3").Select()
